$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells for team record columns
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the header style (bold, bordered, centered) used by the rest of row 1
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Fill in the team's win/loss/tie record for every player row
for ($r = 2; $r -le 51; $r++) {
    $ws.Cells.Item($r, 30).Value = 89   # AD
    $ws.Cells.Item($r, 31).Value = 74   # AE
    $ws.Cells.Item($r, 32).Value = 0    # AF
}
